$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the "Matplotlib" task (old row 5), pushing
# everything from old row 5 onward down by one. Excel carries the date/
# hyperlink formatting down from the row above into the new blank row.
$ws.Rows("5:5").Insert()

# --- Text edits -----------------------------------------------------

# Row 9 (was row 8): tidy up the long RBD/noise comment text (drop the
# extra space before "forces").
$ws.Range("G9").Value2 = "Comments were created to explain the Random Walk process, coefficients and reasoning. For now however the noise component was removed from RBD, noise to be included in the dynamic model (forces and torques)."

# Row 4: "Reading" -> "Done Reading"
$ws.Range("F4").Value2 = "Done Reading"

# Row 8 (was row 7): new status cell "Ongoing " for the TestCases task
$ws.Range("F8").Value2 = "Ongoing "

# --- Column width -----------------------------------------------------
# Give column F an explicit width (matches new Status column usage).
$ws.Columns("F:F").ColumnWidth = 13.666666666666668

# --- View state ---------------------------------------------------------
$ws.Range("F8").Select() | Out-Null
